$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update row 4 & 5 values, then delete row 6 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")

$wsSchedule.Range("B4").Value = 46079.14583333334
$wsSchedule.Range("C4").Value = 7.5
$wsSchedule.Range("D4").Value = 28.35
$wsSchedule.Range("E4").Value = 826.4322495
$wsSchedule.Range("F4").Value = 29.15104936507937
$wsSchedule.Range("A5").Value = 46079.3125
$wsSchedule.Range("B5").Value = 46079.66666666666
$wsSchedule.Range("C5").Value = 8.5
$wsSchedule.Range("D5").Value = 32.13
$wsSchedule.Range("E5").Value = 255.74873025
$wsSchedule.Range("F5").Value = 7.959811087768442

# Row 6 is removed entirely (shifts dimension from A1:F6 to A1:F5)
$wsSchedule.Rows.Item(6).Delete()

# --- Sheet "Detailed": update Price (B), Type (C), Pump_Status (E) for rows 38-97 ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B38").Value = 69.26627999999999
$wsDetailed.Range("B39").Value = 70.36225
$wsDetailed.Range("B40").Value = 71.02005
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 70.36225
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 69.03394
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 57.31
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 57.06
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 37.89
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("B46").Value = 56.98
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 82.42008
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 57.35
$wsDetailed.Range("C49").Value = "historical"
$wsDetailed.Range("B50").Value = 57.06008
$wsDetailed.Range("E50").Value = "ON"
$wsDetailed.Range("B52").Value = 56.98
$wsDetailed.Range("B53").Value = 51.34198
$wsDetailed.Range("B54").Value = 49.63623
$wsDetailed.Range("B55").Value = 48.68524
$wsDetailed.Range("B56").Value = 51.7551
$wsDetailed.Range("B57").Value = 56.98
$wsDetailed.Range("E57").Value = "OFF"
$wsDetailed.Range("B58").Value = 53.45754
$wsDetailed.Range("E58").Value = "OFF"
$wsDetailed.Range("B59").Value = 56.98
$wsDetailed.Range("E59").Value = "OFF"
$wsDetailed.Range("B60").Value = 57.06006
$wsDetailed.Range("B62").Value = 71.39019
$wsDetailed.Range("B63").Value = 76.22794
$wsDetailed.Range("B64").Value = 64.99988
$wsDetailed.Range("E65").Value = "ON"
$wsDetailed.Range("B67").Value = 35.88
$wsDetailed.Range("B68").Value = 13.12638
$wsDetailed.Range("B69").Value = 10.4403
$wsDetailed.Range("B70").Value = 5.91519
$wsDetailed.Range("B71").Value = 0.70613
$wsDetailed.Range("B72").Value = 0.0112
$wsDetailed.Range("B73").Value = 0.5101
$wsDetailed.Range("B74").Value = 0.66949
$wsDetailed.Range("B75").Value = 0.67367
$wsDetailed.Range("B76").Value = 2.45449
$wsDetailed.Range("B77").Value = 0.01078
$wsDetailed.Range("B78").Value = 0.51
$wsDetailed.Range("B79").Value = 37.89
$wsDetailed.Range("B80").Value = 43.8586
$wsDetailed.Range("B81").Value = 35.88
$wsDetailed.Range("B82").Value = 35.88
$wsDetailed.Range("B84").Value = 47.85084
$wsDetailed.Range("B85").Value = 47.7311
$wsDetailed.Range("B86").Value = 47.6007
$wsDetailed.Range("B87").Value = 62.41151
$wsDetailed.Range("B88").Value = 66.70088
$wsDetailed.Range("B89").Value = 71.40000000000001
$wsDetailed.Range("B90").Value = 78
$wsDetailed.Range("B91").Value = 73.2
$wsDetailed.Range("B92").Value = 71.40000000000001
$wsDetailed.Range("B93").Value = 66.10442
$wsDetailed.Range("B94").Value = 64.92106
$wsDetailed.Range("B95").Value = 57.14733
$wsDetailed.Range("B96").Value = 57.06
$wsDetailed.Range("B97").Value = 57.06
